# Update the three-digit-by-one-digit multiplication answers in the table
# to match the regenerated answer key (commit 9a8706d).

$d = $word.ActiveDocument

$replacements = @(
    @("982×5=4910", "818×4=3272"),
    @("341×2=682",  "232×6=1392"),
    @("290×4=1160", "362×3=1086"),
    @("340×2=680",  "610×7=4270"),
    @("943×6=5658", "636×8=5088"),
    @("504×5=2520", "512×6=3072"),
    @("973×2=1946", "389×4=1556"),
    @("673×9=6057", "296×6=1776"),
    @("863×7=6041", "397×4=1588"),
    @("794×9=7146", "626×3=1878"),
    @("841×3=2523", "903×8=7224"),
    @("753×9=6777", "963×3=2889"),
    @("317×4=1268", "110×8=880"),
    @("696×7=4872", "370×4=1480"),
    @("611×8=4888", "561×5=2805"),
    @("512×5=2560", "384×4=1536"),
    @("486×6=2916", "442×9=3978"),
    @("448×9=4032", "253×6=1518"),
    @("792×9=7128", "130×9=1170"),
    @("148×3=444",  "393×9=3537"),
    @("241×4=964",  "186×6=1116"),
    @("365×9=3285", "348×9=3132"),
    @("793×4=3172", "328×9=2952"),
    @("870×8=6960", "601×2=1202"),
    @("217×7=1519", "690×4=2760")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: could not find '$old' to replace with '$new'"
    }
}

Write-Output "Replacements complete"
